# Payroll example layout — "Cambios dispersion separacion bancomer, bug al editar empleado"
#
# The three "NUEVO SUELDO" amounts in the IPSNet Salario sheet were being
# stored/recalculated as raw (sometimes repeating-decimal) numbers, e.g.
# 9666.6666666666679, which caused the bug when editing an employee record
# (Bancomer dispersion/separation). They are replaced with fixed, pre-
# rounded text values so the layout always shows the exact figure that was
# dispersed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPSNet Salario")

# Column G ("NUEVO SUELDO ") goes from a numeric (Millares) format to plain
# text holding the exact amount. Order matters: it controls the order the
# new values land in the shared-string table (G4 first, G3 second, G2
# third; G5 repeats G4's value and reuses that shared string).
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "9666.66 "

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "14666.66 "

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5500.00"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "9666.66 "

# Last selected cell on the sheet moved from D16 to F10.
[void]$ws.Range("F10").Select()
